$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows (56-74) appended to the extraction table, matching the
# source edit: two new studies (garcaruiz2011 / table 1, developmental
# time of Xylotrechus arvicola) plus a trailing study_id/origin stub row
# (glass2019 / figure 2).
# ---------------------------------------------------------------------------
$rowData = @(
"56|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=15;I=15;J=15;K=24;M=""developmental time"";N=""days "";O=29.48;P=0;Q=101;R=0;S=7.0000000000000007E-2;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"57|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=18;I=18;J=18;K=24;M=""developmental time"";N=""days "";O=17.52;P=0;Q=114;R=0;S=0.05;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"58|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=21;I=21;J=21;K=24;M=""developmental time"";N=""days "";O=10.3;P=0;Q=107;R=0;S=7.0000000000000007E-2;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"59|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=24;I=24;J=24;K=24;M=""developmental time"";N=""days "";O=8.02;P=0;Q=115;R=0;S=0.02;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"60|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=27;I=27;J=27;K=24;M=""developmental time"";N=""days "";O=7;P=0;Q=77;R=0;S=0;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"61|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=30;I=30;J=30;K=24;M=""developmental time"";N=""days "";O=6.95;P=0;Q=106;R=0;S=0.03;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"62|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=32;I=32;J=32;K=24;M=""developmental time"";N=""days "";O=6.03;P=0;Q=68;R=0;S=0.09;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"63|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=34;I=34;J=34;K=24;M=""developmental time"";N=""days "";O=6.41;P=0;Q=144;R=0;S=0.09;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"64|A=10;B=""garcaruiz2011"";C=""table 1"";D=0;E=0;G=24;H=35;I=35;J=35;K=24;M=""developmental time"";N=""days "";O=7.48;P=0;Q=178;R=0;S=0.12;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"65|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=15;I=12;J=18.5;K=24;M=""developmental time"";N=""days "";O=29.6;P=0;Q=113;R=0;S=0.32;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"66|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=18;I=12;J=24.5;K=24;M=""developmental time"";N=""days "";O=16.04;P=0;Q=155;R=0;S=0.04;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"67|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=21;I=14.5;J=27.5;K=24;M=""developmental time"";N=""days "";O=11.48;P=0;Q=148;R=0;S=0.08;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"68|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=24;I=17;J=30;K=24;M=""developmental time"";N=""days "";O=10.61;P=0;Q=184;R=0;S=0.04;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"69|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=27;I=17.5;J=30.5;K=24;M=""developmental time"";N=""days "";O=7.56;P=0;Q=186;R=0;S=0.04;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"70|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=30;I=23.5;J=36.5;K=24;M=""developmental time"";N=""days "";O=7.32;P=0;Q=150;R=0;S=7.0000000000000007E-2;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"71|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=32;I=26.5;J=36.5;K=24;M=""developmental time"";N=""days "";O=7.25;P=0;Q=96;R=0;S=0.09;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"72|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=34;I=30.5;J=36.5;K=24;M=""developmental time"";N=""days "";O=7.31;P=0;Q=128;R=0;S=0.15;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"73|A=11;B=""garcaruiz2011"";C=""table 1"";D=0;E=1;F=1;G=24;H=35;I=32;J=36.5;K=24;M=""developmental time"";N=""days "";O=8;P=0;Q=103;R=0;T=""Xylotrechus "";U=""arvicola"";V=1;W=1;X=1",
"74|B=""glass2019"";C=""figure 2"""
)

foreach ($rowEntry in $rowData) {
    $rowParts = $rowEntry -split "\|", 2
    $rowNum = $rowParts[0]
    $cellsStr = $rowParts[1]
    $cellDefs = $cellsStr -split ";"
    foreach ($cellDef in $cellDefs) {
        $eqIdx = $cellDef.IndexOf("=")
        $colLetter = $cellDef.Substring(0, $eqIdx)
        $rawVal = $cellDef.Substring($eqIdx + 1)
        $cellRef = "$colLetter$rowNum"
        if ($rawVal.StartsWith('"')) {
            $strVal = $rawVal.Substring(1, $rawVal.Length - 2)
            $ws.Range($cellRef).Value = $strVal
        } else {
            $numVal = [double]$rawVal
            $ws.Range($cellRef).Value = $numVal
        }
    }
}

# ---------------------------------------------------------------------------
# Carry over the existing cell formatting (fonts) used by the table for the
# resp_quality / samp_size / larger_group / exp_age / size columns so the
# new rows 56-73 match the look of the preceding rows (row 55 is used as
# the formatting template, same as the rest of the sheet from row 8 on).
# ---------------------------------------------------------------------------
$ws.Range("P55").Copy()
$ws.Range("P56:P73").PasteSpecial(-4122)

$ws.Range("Q55").Copy()
$ws.Range("Q56:Q73").PasteSpecial(-4122)

$ws.Range("V55").Copy()
$ws.Range("V56:V73").PasteSpecial(-4122)

$ws.Range("W55").Copy()
$ws.Range("W56:W73").PasteSpecial(-4122)

$ws.Range("X55").Copy()
$ws.Range("X56:X73").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selection / view bookkeeping to mirror the author's saved view state.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 480
$win.Top = 600

$ws.Range("C74").Select()
